$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "68.719.36"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.59%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.712.20"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.38%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "599.78"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "163.11"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.54%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.22%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.711.20"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("E13").Value = "  +2.67%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "28.46"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.30%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.209.95"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("E16").Value = "  -0.27%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "68.636.45"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.49%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.719.18"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.14%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.89"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.34%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.68"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +4.26%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "365.36"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  +2.44%  "
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("E24").Value = "  +2.69%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "73.81"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.91%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  +1.76%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.842.55"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.65%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "594.62"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +6.35%  "
$ws.Range("E31").Value = "  +0.03%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.24"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.19%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.96"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +5.16%  "
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("E36").Value = "  +4.57%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "161.17"
$c.Style = "Normal"
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "19.92"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("E43").Value = "  +2.99%  "
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("E46").Value = "  -5.57%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "157.87"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("E48").Value = "  +5.12%  "
$ws.Range("E49").Value = "  +5.87%  "
$ws.Range("E50").Value = "  +7.07%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "22.04"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.48%  "
